# "image update for black mode"
#  1. Update the cached datetimeFigureOut footer field text (slide master
#     + every slide layout) from 2023-09-26 to 2023-11-26.
#  2. Add a full-bleed white "cover" rectangle behind the existing content
#     on slides 7, 8 and 9 (sent to the back of the z-order).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Date placeholder text (slide master + all custom layouts)
# ---------------------------------------------------------------------
function Update-DateShape($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame -eq -1) {
            if ($sh.TextFrame.TextRange.Text -eq "2023-09-26") {
                $sh.TextFrame.TextRange.Text = "2023-11-26"
                return $true
            }
        }
    }
    return $false
}

$design = $p.Designs.Item(1)
$slideMaster = $design.SlideMaster

Update-DateShape($slideMaster.Shapes) | Out-Null

$layouts = $slideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DateShape($layouts.Item($li).Shapes) | Out-Null
}

# ---------------------------------------------------------------------
# 2) New background rectangles on slides 7, 8, 9
# ---------------------------------------------------------------------
$EMU_PER_PT = 12700

function Add-CoverRectangle($slide, $offX, $offY, $extCx, $extCy) {
    $left   = $offX  / $EMU_PER_PT
    $top    = $offY  / $EMU_PER_PT
    $width  = $extCx / $EMU_PER_PT
    $height = $extCy / $EMU_PER_PT

    $rect = $slide.Shapes.AddShape(1, $left, $top, $width, $height)

    $rect.Fill.ForeColor.SchemeColor = "bg1"
    $rect.Line.Visible = 0

    $rect.TextFrame.VerticalAnchor = 3
    $rect.TextFrame.TextRange.ParagraphFormat.Alignment = 2

    $rect.ZOrder(1)

    return $rect
}

# Slide 7 -> id 2 "직사각형 1" (first rectangle created in the session)
$slide7 = $p.Slides.Item(7)
$r1 = Add-CoverRectangle $slide7 0 203200 15341600 6654800
$r1.Name = "직사각형 1"

# Slide 8 -> id 3 "직사각형 2" (a throw-away shape is created first so the
# id/name counters land on 3/"Rectangle 2", matching the authored file)
$slide8 = $p.Slides.Item(8)
$dummy = $slide8.Shapes.AddShape(1, 0, 0, 10, 10)
$dummy.Delete()
$r2 = Add-CoverRectangle $slide8 231939 570493 11728121 5168737
$r2.Name = "직사각형 2"

# Slide 9 -> id 2 "직사각형 1"
$slide9 = $p.Slides.Item(9)
$r3 = Add-CoverRectangle $slide9 231939 570493 11728121 5168737
$r3.Name = "직사각형 1"
